$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update data values
$ws.Range("B2").Value = 6000
$ws.Range("B3").Value = 10500

# Update the selected cell on the sheet
$ws.Range("B3").Select()
